$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.065.22"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "1.783.31"
$ws.Range("E3").Value = "  -2.50%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "

$ws.Range("E6").Value = "  -1.70%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.49"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.38%  "

$ws.Range("E9").Value = "  -2.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0710"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.56%  "

$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").Value = "2.041.93"
$ws.Range("E12").Value = "  -2.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.97%  "

$ws.Range("D14").Value = "1.794.91"
$ws.Range("E14").Value = "  -1.97%  "

$ws.Range("D15").Value = "33.997.20"
$ws.Range("E15").Value = "  -0.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.618"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.16%  "

$ws.Range("D20").Value = "0.0₃0785"
$ws.Range("E20").Value = "  -1.16%  "

$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.77%  "

$ws.Range("E28").Value = "  -2.71%  "

$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("E30").Value = "  +0.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0512"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.40%  "

$ws.Range("E32").Value = "  -4.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.31%  "

$ws.Range("D35").Value = "1.390.41"
$ws.Range("E35").Value = "  -3.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.644"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.04"
$ws.Range("D37").Style = "Normal"

$ws.Range("E38").Value = "  -1.99%  "

$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.10%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.910"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.10%  "

$ws.Range("E42").Value = "  -2.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "77.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.07%  "

$ws.Range("E44").Value = "  +13.77%  "

$ws.Range("E45").Value = "  +2.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0497"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.19%  "

$ws.Range("E49").Value = "  -4.54%  "

$ws.Range("D50").Value = "1.939.77"
$ws.Range("E50").Value = "  -2.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.19%  "
